# "Add files via upload" — make the valency-class / X / Y columns uniform:
# for every data row where column I ("X") currently holds the TR locus
# label, split it into I="ERG" (unchanged pattern, just re-set) and give
# column J ("Y") the matching ABS label (it was previously left blank).
# A handful of rows instead carry the "*" placeholder in I; for those, J
# gets the same "*" placeholder instead of ABS.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Udi")

# Rows whose column I is "TR" -> becomes "ERG", column J becomes "ABS"
$ergRows = @(5,9,10,11,14,16,17,21,23,27,28,30,32,34,35,37,40,41,42,44,45,46,47,50,51,53,56,58,61,64,67,70,71,72,73,76,79,80,86,87,92,94,97,98,102,103,104,106,107,108,109,110,111,114,115,116,120,122,125)

foreach ($r in $ergRows) {
    $ws.Cells.Item($r, 9).Value = "ERG"
    $ws.Cells.Item($r, 10).Value = "ABS"
}

# Rows whose column I is the placeholder "*" -> stays "*", column J becomes "*" too
$starRows = @(49,81,84,93,121,126,127)

foreach ($r in $starRows) {
    $ws.Cells.Item($r, 9).Value = "*"
    $ws.Cells.Item($r, 10).Value = "*"
}

# Reset the saved viewport/selection state (the sheet previously had
# topLeftCell="I1" / a selection parked at O3 - normalise it back).
$aw = $excel.ActiveWindow
$aw.ScrollColumn = 1
$aw.ScrollRow = 1
$ws.Range("A1").Select()
